$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.713.71"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.28%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'1.634.49"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'213.73"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.62%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'0.499"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.09%  "
$ws.Range("E6").ClearFormats()
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("E7").ClearFormats()
$ws.Range("E8").Value = "  +1.34%  "
$ws.Range("E8").ClearFormats()
$ws.Range("E9").Value = "  +0.98%  "
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'19.25"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.10%  "
$ws.Range("E10").ClearFormats()
$ws.Range("E11").Value = "  +3.39%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'1.861.79"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.63%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'1.670.79"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.97%  "
$ws.Range("E13").ClearFormats()
$ws.Range("E14").Value = "  +1.72%  "
$ws.Range("E14").ClearFormats()
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'26.708.86"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.24%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'63.62"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.66%  "
$ws.Range("E17").ClearFormats()
$ws.Range("E18").Value = "  +2.10%  "
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'218.27"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +7.63%  "
$ws.Range("E19").ClearFormats()
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("E20").ClearFormats()
$ws.Range("E21").Value = "  +1.17%  "
$ws.Range("E21").ClearFormats()
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'6.15"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.50%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'1.96"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +5.12%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'147.80"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.07%  "
$ws.Range("E25").ClearFormats()
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("E26").ClearFormats()
$ws.Range("E27").Value = "  +1.34%  "
$ws.Range("E27").ClearFormats()
$ws.Range("E28").Value = "  +3.99%  "
$ws.Range("E28").ClearFormats()
$ws.Range("E29").Value = "  +2.42%  "
$ws.Range("E29").ClearFormats()
$ws.Range("E30").Value = "  -2.16%  "
$ws.Range("E30").ClearFormats()
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("E31").ClearFormats()
$ws.Range("D32").Value = "'3.32"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.73%  "
$ws.Range("E32").ClearFormats()
$ws.Range("E33").Value = "  +2.41%  "
$ws.Range("E33").ClearFormats()
$ws.Range("E34").Value = "  +0.61%  "
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = "'1.228.88"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +5.88%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'2.39"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E36").ClearFormats()
$ws.Range("E37").Value = "  +5.78%  "
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'0.807"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.47%  "
$ws.Range("E38").ClearFormats()
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'0.501"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.16%  "
$ws.Range("E40").ClearFormats()
$ws.Range("D42").Value = "'0.797"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.95%  "
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'5.35"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.71%  "
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'1.770.87"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'93.00"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.84%  "
$ws.Range("E45").ClearFormats()
$ws.Range("E46").Value = "  +2.82%  "
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'55.48"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.76%  "
$ws.Range("E47").ClearFormats()
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("E48").ClearFormats()
$ws.Range("E49").Value = "  +1.20%  "
$ws.Range("E49").ClearFormats()
$ws.Range("E50").Value = "  +4.48%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'0.408"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.17%  "
$ws.Range("E51").ClearFormats()
